$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "PSO optimisation "
$ws.Range("C9").Value = "Random"
$ws.Range("D9").Value = "NA"

$ws.Range("D10").Select()
